# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values on the "zh-cn" and "de-de" sheets to reflect a freshly regenerated
# report. The cells hold plain text (not real Excel dates), so set them as
# strings with Value2/Value to avoid Excel converting them into date serials.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 08:19:05"
$wsZhCn.Range("E5").Value = "2016-03-19 08:19:05"
$wsZhCn.Range("H2").Value = "2016-03-19 08:19:27"
$wsZhCn.Range("H5").Value = "2016-03-19 08:19:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 08:19:08"
$wsDeDe.Range("E5").Value = "2016-03-19 08:19:08"
$wsDeDe.Range("H2").Value = "2016-03-19 08:19:32"
$wsDeDe.Range("H5").Value = "2016-03-19 08:19:32"
